$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

for ($row = 2; $row -le 12; $row++) {
    # Column A ("card") holds a number-looking value that must stay text,
    # same as the original "2" -> leading apostrophe forces text storage.
    $ws.Cells.Item($row, 1).Value = "'22"

    # Column O ("Servised by") was blank; fill with "nan" like the rest
    # of the sheet's empty cells.
    $ws.Cells.Item($row, 15).Value = "nan"
}
